$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(42625.885810185187)
$ws.Cells.Item($row, 2).Value = 30
$ws.Cells.Item($row, 3).Value = 69
$ws.Cells.Item($row, 4).Value = 30
$ws.Cells.Item($row, 5).Value = 99
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 11675
$ws.Cells.Item($row, 8).Value = 6361
$ws.Cells.Item($row, 9).Value = 326
$ws.Cells.Item($row, 10).Value = 108
$ws.Cells.Item($row, 11).Value = 47
$ws.Cells.Item($row, 12).Value = 7
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Named"
